$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "BBB"
$ws.Range("B9").Value = "A"
$ws.Range("B10").Value = "A"
$ws.Range("B11").Value = "A"
$ws.Range("B12").Value = "AAA"
$ws.Range("B13").Value = "BBB"
$ws.Range("B18").Value = "BBB"
$ws.Range("B19").Value = "BB"
$ws.Range("B22").Value = "B"
$ws.Range("B23").Value = "B"
$ws.Range("B24").Value = "A"
$ws.Range("B28").Value = "A"
$ws.Range("B29").Value = "BB"
$ws.Range("B30").Value = "AA"
$ws.Range("B32").Value = "A"
$ws.Range("B33").Value = "BBB"
$ws.Range("B37").Value = "BBB"
$ws.Range("B42").Value = "BBB"
$ws.Range("B43").Value = "BBB"
$ws.Range("B46").Value = "BB"
$ws.Range("B50").Value = "A"
$ws.Range("B51").Value = "BB"
$ws.Range("B52").Value = "BB"
$ws.Range("B55").Value = "BB"
$ws.Range("B57").Value = "A"
$ws.Range("B59").Value = "BB"
$ws.Range("B60").Value = "BB"
$ws.Range("B65").Value = "BB"
$ws.Range("B66").Value = "BBB"
$ws.Range("B69").Value = "BB"
$ws.Range("B73").Value = "A"
$ws.Range("B74").Value = "BB"
$ws.Range("B76").Value = "BBB"
$ws.Range("B83").Value = "AA"
$ws.Range("B85").Value = "BB"
$ws.Range("B90").Value = "BBB"
$ws.Range("B91").Value = "A"
$ws.Range("B92").Value = "BBB"
$ws.Range("B102").Value = "BBB"
$ws.Range("B106").Value = "AA"
$ws.Range("B110").Value = "BBB"
$ws.Range("B111").Value = "BBB"
$ws.Range("B116").Value = "BBB"
$ws.Range("B117").Value = "BBB"
$ws.Range("B119").Value = "BBB"
$ws.Range("B120").Value = "BBB"
$ws.Range("B124").Value = "BBB"
$ws.Range("B129").Value = "BBB"
$ws.Range("B134").Value = "A"
$ws.Range("B135").Value = "BBB"
$ws.Range("B138").Value = "BBB"
$ws.Range("B142").Value = "A"
$ws.Range("B144").Value = "A"
$ws.Range("B145").Value = "A"
$ws.Range("B154").Value = "BBB"
$ws.Range("B156").Value = "B"
$ws.Range("B158").Value = "BB"
$ws.Range("B160").Value = "A"
$ws.Range("B164").Value = "B"
$ws.Range("B165").Value = "B"
$ws.Range("B166").Value = "BBB"
$ws.Range("B167").Value = "BBB"
$ws.Range("B168").Value = "BB"
$ws.Range("B170").Value = "BBB"
$ws.Range("B173").Value = "B"
$ws.Range("B177").Value = "A"
$ws.Range("B178").Value = "A"
$ws.Range("B180").Value = "AAA"
$ws.Range("B181").Value = "AA"
$ws.Range("B182").Value = "BBB"
$ws.Range("B185").Value = "BB"
$ws.Range("B187").Value = "AAA"
$ws.Range("B189").Value = "BBB"
$ws.Range("B194").Value = "BBB"
$ws.Range("B199").Value = "BB"
$ws.Range("B200").Value = "BB"
$ws.Range("B203").Value = "A"
$ws.Range("B204").Value = "AAA"
$ws.Range("B205").Value = "AAA"
$ws.Range("B206").Value = "BB"
$ws.Range("B209").Value = "BBB"
$ws.Range("B212").Value = "A"
$ws.Range("B215").Value = "BB"
$ws.Range("B224").Value = "BBB"
$ws.Range("B227").Value = "BBB"
$ws.Range("B228").Value = "B"
$ws.Range("B237").Value = "BB"
$ws.Range("B239").Value = "AAA"
$ws.Range("B240").Value = "AAA"
$ws.Range("B241").Value = "AAA"
$ws.Range("B245").Value = "BBB"
$ws.Range("B252").Value = "BB"
$ws.Range("B255").Value = "BB"
$ws.Range("B256").Value = "BBB"
$ws.Range("B257").Value = "BB"
$ws.Range("B259").Value = "BBB"
$ws.Range("B262").Value = "A"
$ws.Range("B272").Value = "BBB"
$ws.Range("B274").Value = "BB"
$ws.Range("B275").Value = "BBB"
$ws.Range("B278").Value = "A"
$ws.Range("B279").Value = "A"
$ws.Range("B280").Value = "A"
$ws.Range("B282").Value = "A"
$ws.Range("B286").Value = "BBB"
$ws.Range("B287").Value = "BBB"
$ws.Range("B288").Value = "A"
$ws.Range("B290").Value = "BBB"
$ws.Range("B291").Value = "A"
$ws.Range("B292").Value = "A"
$ws.Range("B293").Value = "A"
$ws.Range("B302").Value = "AA"
$ws.Range("B303").Value = "AAA"
$ws.Range("B304").Value = "BB"
$ws.Range("B305").Value = "BBB"
$ws.Range("B306").Value = "BBB"
$ws.Range("B312").Value = "BBB"
$ws.Range("B314").Value = "BB"
$ws.Range("B316").Value = "BB"
$ws.Range("B319").Value = "BBB"
$ws.Range("B320").Value = "BBB"
$ws.Range("B321").Value = "BBB"
$ws.Range("B322").Value = "BBB"
$ws.Range("B330").Value = "A"
$ws.Range("B331").Value = "BBB"
$ws.Range("B332").Value = "B"
$ws.Range("B334").Value = "BB"
$ws.Range("B335").Value = "BB"
$ws.Range("B342").Value = "BB"
$ws.Range("B343").Value = "BBB"
$ws.Range("B347").Value = "BBB"
$ws.Range("B348").Value = "BB"
$ws.Range("B350").Value = "B"
$ws.Range("B357").Value = "BBB"
$ws.Range("B358").Value = "BBB"
$ws.Range("B360").Value = "AA"
$ws.Range("B362").Value = "AA"
$ws.Range("B366").Value = "A"
$ws.Range("B368").Value = "BBB"
$ws.Range("B370").Value = "BBB"
$ws.Range("B376").Value = "AA"
$ws.Range("B382").Value = "B"
$ws.Range("B384").Value = "BB"
$ws.Range("B387").Value = "B"
$ws.Range("B388").Value = "BB"
$ws.Range("B396").Value = "B"
$ws.Range("B397").Value = "BBB"
$ws.Range("B401").Value = "AAA"
